# add: new function codelco
# Normalize column A (Pos) to "00010" and column B (Material) to "11111"
# for every data row, renumber column E (Bulto) sequentially, and fix up
# a handful of column C (Cantidad) values to match the new target data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2..18: Pos(A), Material(B), Cantidad(C), Unidad(D), Bulto(E)
$data = @(
    @("00010", "11111", 4, "UN", 1),
    @("00010", "11111", 4, "UN", 1),
    @("00010", "11111", 4, "UN", 2),
    @("00010", "11111", 5, "UN", 3),
    @("00010", "11111", 4, "UN", 4),
    @("00010", "11111", 5, "UN", 5),
    @("00010", "11111", 4, "UN", 6),
    @("00010", "11111", 4, "UN", 7),
    @("00010", "11111", 5, "UN", 8),
    @("00010", "11111", 4, "UN", 9),
    @("00010", "11111", 4, "UN", 10),
    @("00010", "11111", 4, "UN", 11),
    @("00010", "11111", 4, "UN", 12),
    @("00010", "11111", 4, "UN", 13),
    @("00010", "11111", 4, "UN", 14),
    @("00010", "11111", 4, "UN", 15),
    @("00010", "11111", 3, "UN", 16)
)

$startRow = 2
$endRow = $startRow + $data.Count - 1

# Columns A and B hold text-like values ("00010", "11111") which must stay
# text (not be coerced into numbers), so force text formatting first.
$ws.Range("A$startRow" + ":A$endRow").NumberFormat = "@"
$ws.Range("B$startRow" + ":B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
